{"js": "// Abstract update: the combined Affera effect estimate and its CI changed.\n//   \"...was a 15.0% reduction (95% CI 8.2\u201321.3%).\"\n// becomes\n//   \"...was a 15.1% reduction (95% CI 10.6\u201319.4%).\"\nconst body = context.document.body;\n\nconst oldClause = \"was a 15.0% reduction (95% CI 8.2\\u201321.3%).\";\nconst newClause = \"was a 15.1% reduction (95% CI 10.6\\u201319.4%).\";\n\nconst results = body.search(oldClause, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target sentence for the Affera combined-effect statistic was not found.\");\n}\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newClause, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Abstract update: the combined Affera effect estimate and its CI changed.\n#   \"...was a 15.0% reduction (95% CI 8.2\u201321.3%).\"\n# becomes\n#   \"...was a 15.1% reduction (95% CI 10.6\u201319.4%).\"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"was a 15.0% reduction (95% CI 8.2\u201321.3%).\"\n$find.Replacement.Text = \"was a 15.1% reduction (95% CI 10.6\u201319.4%).\"\n$find.Forward = $true\n$find.Wrap = 0\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n\n$found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 0, $false, $find.Replacement.Text, 2)\n\nif (-not $found) {\n  throw \"Target sentence for the Affera combined-effect statistic was not found.\"\n}\n\nWrite-Output \"done\"\n"}
